# Update countries & provincias Spain
# Refreshes the COVID-19 "Pais" dashboard sheet with newer source data:
#  - Updates the "last refreshed" timestamp.
#  - Updates case/death counters for several countries whose source figures
#    changed between pulls.
#  - Three pairs of neighbouring countries swapped places in the sort order
#    (the sheet is sorted by total cases, column B, descending), so both the
#    country name (column A) and its row of statistics move together.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# --- "Datos actualizados" footer timestamp -------------------------------
$ws.Range("A1").Value = "Datos actualizados a 4 de Octubre de 2020 a las 22:46"

# --- Estados Unidos (row 4) ------------------------------------------------
$ws.Range("B4").Value = 7628637
$ws.Range("C4").Value = 26785
$ws.Range("D4").Value = 4840048
$ws.Range("E4").Value = 2574074
$ws.Range("G4").Value = 236
$ws.Range("H4").Value = 214515

# --- Sudafrica (row 13) ----------------------------------------------------
$ws.Range("B13").Value = 681289
$ws.Range("C13").Value = 1573
$ws.Range("D13").Value = 614781
$ws.Range("E13").Value = 49532
$ws.Range("G13").Value = 38
$ws.Range("H13").Value = 16976

# --- Alemania (row 26) ------------------------------------------------------
$ws.Range("B26").Value = 301543
$ws.Range("C26").Value = 1515
$ws.Range("E26").Value = 30041
$ws.Range("G26").Value = 5
$ws.Range("H26").Value = 9602

# --- Israel (row 27) --------------------------------------------------------
$ws.Range("B27").Value = 266775
$ws.Range("C27").Value = 2332
$ws.Range("D27").Value = 195629
$ws.Range("E27").Value = 69427
$ws.Range("G27").Value = 37
$ws.Range("H27").Value = 1719

# --- Angola (row 124) --------------------------------------------------------
$ws.Range("B124").Value = 5402
$ws.Range("C124").Value = 32
$ws.Range("D124").Value = 2577
$ws.Range("E124").Value = 2630
$ws.Range("G124").Value = 2
$ws.Range("H124").Value = 195

# --- Togo / Nueva Zelanda swap places (rows 160-161) ------------------------
$ws.Range("A160").Value = "Togo"
$ws.Range("C160").Value = 14
$ws.Range("D160").Value = 1392
$ws.Range("E160").Value = 414
$ws.Range("H160").Value = 48

$ws.Range("A161").Value = "Nueva Zelanda"
$ws.Range("B161").Value = 1854
$ws.Range("C161").Value = 5
$ws.Range("D161").Value = 1788
$ws.Range("E161").Value = 41
$ws.Range("H161").Value = 25

# --- Curazao / Gibraltar swap places (rows 180-181) -------------------------
$ws.Range("A180").Value = "Curazao"
$ws.Range("B180").Value = 452
$ws.Range("C180").Value = 23
$ws.Range("D180").Value = 231
$ws.Range("E180").Value = 220
$ws.Range("H180").Value = 1

$ws.Range("A181").Value = "Gibraltar"
$ws.Range("B181").Value = 432
$ws.Range("C181").Value = 4
$ws.Range("D181").Value = 360
$ws.Range("E181").Value = 72
$ws.Range("H181").Value = 0

# --- Monaco (row 189) --------------------------------------------------------
$ws.Range("D189").Value = 189
$ws.Range("E189").Value = 31

# --- Montserrat / Islas Malvinas swap places (rows 215-216) -----------------
$ws.Range("A215").Value = "Montserrat"
$ws.Range("D215").Value = 12
$ws.Range("H215").Value = 1

$ws.Range("A216").Value = "Islas Malvinas"
$ws.Range("D216").Value = 13
$ws.Range("H216").Value = 0
